# UNI-0002 <I> First big update; updated button images, logic controller,
# LHDataLevels, LHUIButton extended from UIButton
#
# Adds a new ticket row to the "Main" sheet of the bug tracker. Everything
# downstream (Summary sheet COUNTIF totals, chart caches) recalculates
# automatically off of this new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Cells.Item(3, 1).Value = "Feature"
$ws.Cells.Item(3, 2).Value = "UNI-0002"
$ws.Cells.Item(3, 3).Value = "<I> First big update. New image for buttons, updated controller logic, updated LHDataLevels, extended UIButton to LHUIButton"
$ws.Cells.Item(3, 4).Value = "rhdelaro"
$ws.Cells.Item(3, 5).Value = "rhdelaro"
$ws.Cells.Item(3, 6).Value = "Low"
$ws.Cells.Item(3, 7).Value = "Development"
$ws.Cells.Item(3, 8).Value = "Unresolved"
$ws.Cells.Item(3, 9).Value = "11/18/2013 06:40:18"

# The summary text wraps onto two lines in the (wrap-text) column C style,
# so the row grows from the default 14pt to 28pt.
$ws.Rows.Item(3).RowHeight = 28
